$d = $word.ActiveDocument
$t = $d.Tables(1)

$values = @(
    "15+29=",
    "19+68=",
    "19+19=",
    "29+0=",
    "14+40=",
    "85-5=",
    "99-80=",
    "51-34=",
    "84-32=",
    "76-30=",
    "97-75=",
    "94-74=",
    "55-16=",
    "13+8=",
    "18+56=",
    "34+36=",
    "83-40=",
    "48+43=",
    "15-0=",
    "10+2=",
    "5+72=",
    "40+27=",
    "94+4=",
    "84-41=",
    "41-18=",
    "63-59=",
    "32+23=",
    "22+56=",
    "62-32=",
    "65-20=",
    "46+46=",
    "57+25=",
    "19+62=",
    "1+67=",
    "21+17=",
    "6+61=",
    "85+5=",
    "39+38=",
    "47+15=",
    "49-45=",
    "57-15=",
    "40+3=",
    "97-7=",
    "73-69=",
    "87-10=",
    "89-45=",
    "44-27=",
    "82-10=",
    "23+49=",
    "71-45=",
    "3+84=",
    "25-18=",
    "88-69=",
    "55-36=",
    "65+30=",
    "2+84=",
    "35-16=",
    "23+18=",
    "49+34=",
    "80-48=",
    "38+50=",
    "76-20=",
    "42-33=",
    "36+17=",
    "97-42=",
    "76-1=",
    "50+31=",
    "20+19=",
    "77+7=",
    "72-17=",
    "56+18=",
    "86+1=",
    "67-32=",
    "41+49=",
    "54+32=",
    "71-67=",
    "24-7=",
    "3+32=",
    "26+7=",
    "44+4=",
    "56-51=",
    "21-20=",
    "91-50=",
    "95-24=",
    "38-13=",
    "36+31=",
    "87-2=",
    "11-5=",
    "96-94=",
    "78-47=",
    "72-54=",
    "88-0=",
    "85-1=",
    "82-28=",
    "36+34=",
    "90-60=",
    "19+62=",
    "18+26=",
    "39+52=",
    "28+71="
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = [math]::Floor($i / 5) + 1
    $col = ($i % 5) + 1
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $values[$i]
}

Write-Host "Done updating" $values.Length "cells"